$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5899
$ws.Range("J62").Value = 5899
$ws.Range("L62").Value = 5899
$ws.Range("N62").Value = -7147
$ws.Range("H65").Value = 5899
$ws.Range("J65").Value = 5899
$ws.Range("L65").Value = 29495
$ws.Range("N65").Value = -35735
$ws.Range("H70").Value = 800
$ws.Range("I70").Value = 800
$ws.Range("K70").Value = 2400
$ws.Range("M70").Value = -2130
$ws.Range("H73").Value = 800
$ws.Range("I73").Value = 800
$ws.Range("K73").Value = 2400
$ws.Range("M73").Value = -1464
$ws.Range("H112").Value = 1027.3334
$ws.Range("J112").Value = 1027.3334
$ws.Range("L112").Value = 3082.0002
$ws.Range("N112").Value = -5298.0002
$ws.Range("H138").Value = 9160.947
$ws.Range("J138").Value = 9797.647000000001
$ws.Range("L138").Value = 29392.941
$ws.Range("N138").Value = -39672.94100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H45").Value = 1989.5
$ws.Range("I45").Value = 1989.5
$ws.Range("K45").Value = 1989.5
$ws.Range("M45").Value = -1612.5
$ws.Range("H61").Value = 6318.4287
$ws.Range("I61").Value = 6816
$ws.Range("J61").Value = 3333
$ws.Range("K61").Value = 6816
$ws.Range("L61").Value = 3333
$ws.Range("M61").Value = -6604
$ws.Range("N61").Value = -3757
$ws.Range("H74").Value = 738.6667
$ws.Range("I74").Value = 722.1429000000001
$ws.Range("K74").Value = 722.1429000000001
$ws.Range("M74").Value = 151.8570999999999
$ws.Range("H77").Value = 738.6667
$ws.Range("I77").Value = 722.1429000000001
$ws.Range("K77").Value = 3610.7145
$ws.Range("M77").Value = 757.2855
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 3261.3125
$ws.Range("I132").Value = 2638.6
$ws.Range("K132").Value = 7915.799999999999
$ws.Range("M132").Value = -5385.799999999999
$ws.Range("H136").Value = 6318.4287
$ws.Range("I136").Value = 6816
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 20448
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -17898
$ws.Range("N136").Value = -15099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H20").Value = 1852.2858
$ws.Range("I20").Value = 1775.8
$ws.Range("J20").Value = 2043.5
$ws.Range("K20").Value = 1775.8
$ws.Range("L20").Value = 2043.5
$ws.Range("M20").Value = -1528.8
$ws.Range("N20").Value = -2537.5
$ws.Range("H86").Value = 3502.2856
$ws.Range("I86").Value = 3502.2856
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3502.2856
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2379.2856
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3502.2856
$ws.Range("I89").Value = 3502.2856
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17511.428
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11895.428
$ws.Range("N89").ClearContents()
$ws.Range("H134").Value = 3838.6428
$ws.Range("I134").Value = 3476.6667
$ws.Range("K134").Value = 10430.0001
$ws.Range("M134").Value = -7895.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 34988.332
$ws.Range("J68").Value = 34988.332
$ws.Range("L68").Value = 34988.332
$ws.Range("N68").Value = -36486.332
$ws.Range("H71").Value = 34988.332
$ws.Range("J71").Value = 34988.332
$ws.Range("L71").Value = 104964.996
$ws.Range("N71").Value = -112452.996
$ws.Range("H132").Value = 4781.579
$ws.Range("I132").Value = 4457.067
$ws.Range("K132").Value = 13371.201
$ws.Range("M132").Value = -10841.201
$ws.Range("H134").Value = 3518.9
$ws.Range("I134").Value = 3493.5789
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 10480.7367
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -7945.736699999999
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7374.5
$ws.Range("I70").Value = 6999
$ws.Range("J70").Value = 7750
$ws.Range("K70").Value = 6999
$ws.Range("L70").Value = 7750
$ws.Range("M70").Value = -6729
$ws.Range("N70").Value = -8290
$ws.Range("H73").Value = 7374.5
$ws.Range("I73").Value = 6999
$ws.Range("J73").Value = 7750
$ws.Range("K73").Value = 6999
$ws.Range("L73").Value = 7750
$ws.Range("M73").Value = -6063
$ws.Range("N73").Value = -9622
$ws.Range("H80").Value = 6065.25
$ws.Range("J80").Value = 10005.5
$ws.Range("L80").Value = 10005.5
$ws.Range("N80").Value = -12001.5
$ws.Range("H83").Value = 6065.25
$ws.Range("J83").Value = 10005.5
$ws.Range("L83").Value = 50027.5
$ws.Range("N83").Value = -60011.5
$ws.Range("H102").Value = 3313.9
$ws.Range("I102").Value = 3313.9
$ws.Range("K102").Value = 3313.9
$ws.Range("M102").Value = -1691.9
$ws.Range("H132").Value = 2898.1667
$ws.Range("I132").Value = 2506.75
$ws.Range("K132").Value = 7520.25
$ws.Range("M132").Value = -4990.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 918999.7
$ws.Range("J43").Value = 918999.7
$ws.Range("L43").Value = 918999.7
$ws.Range("N43").Value = -919385.7
$ws.Range("H132").Value = 2408.95
$ws.Range("I132").Value = 1479
$ws.Range("K132").Value = 4437
$ws.Range("M132").Value = -1907

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 15000
$ws.Range("K62").Value = 15000
$ws.Range("M62").Value = -14376
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 15000
$ws.Range("K65").Value = 75000
$ws.Range("M65").Value = -71880
$ws.Range("H136").Value = 1999.5
$ws.Range("I136").Value = 1999.5
$ws.Range("K136").Value = 5998.5
$ws.Range("M136").Value = -3448.5
